$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (the "Förändrad" column) holds a date serial number that was
# bumped by one day (46081 -> 46082) for every data row (rows 2 through 152).
$ws.Range("C2:C152").Value2 = 46082
